$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Rows.Item(2).Insert()
$ws.Range("A4:C4").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
$ws.Range("A2").Value = 45657
$ws.Range("B2").Value = 10.7
$ws.Range("C2").Formula = "=(B2/B3-1)*100"
"done"
